$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing header + data rows down by two rows, opening up a
# blank spacer row (row 2) between the new export-date line and the
# (now relocated) header row on row 3.
$ws.Rows("1:2").Insert()

# Row 1: new "Export Date and Time" stamp, bold like the header row but
# as its own distinct style.
$ws.Cells.Item(1, 1).Value = "Export Date and Time: 2024-05-29 13:12:27"
$ws.Cells.Item(1, 1).Font.Bold = $true
